# CylindricalGauge checklist workbook update
# - Adds two new "Basic Chart" BVT rows (categorical-index checklist entries)
# - Widens column E on the BVT sheet to fit the new, longer text
# - Grows several wrapped-text row heights on the Checklist sheet to match
#   the content after the newer Excel build recalculated wrapping
# - Leaves the active selection on E9, matching the author's final position

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. BVT sheet: two new rows documenting the "Basic Chart" test case
# ---------------------------------------------------------------------
$bvt = $wb.Worksheets.Item("BVT")

# Row 9
$bvt.Range("A9").Value2 = 5
$bvt.Range("B9").Value2 = "Basic Chart"
$bvt.Range("C9").Value2 = "Display basic Cylindrical Gauge"
$bvt.Range("D9").Value2 = "1. Drag 'Revenue' in 'Actual Value' field."
$bvt.Range("E9").Value2 = "Visual should be rendered for Revenue value."

# Row 10
$bvt.Range("A10").Value2 = 6
$bvt.Range("B10").Value2 = "Basic Chart"
$bvt.Range("C10").Value2 = "Display basic Cylindrical Gauge"
$bvt.Range("D10").Value2 = "1. Drag same value in all field."
$bvt.Range("E10").Value2 = "Visual should be rendered properly."

# Match the surrounding rows: content top-aligned (C:E already wrap via the
# column style; only vertical alignment needs to be (re)asserted on the new
# rows since freshly written cells default to bottom alignment).
$bvt.Range("C9:E10").VerticalAlignment = -4160   # xlTop

# ---------------------------------------------------------------------
# 2. BVT sheet: column E needs to be a bit wider to fit the new text and
#    was re-measured (best-fit) by Excel.
# ---------------------------------------------------------------------
$bvt.Columns.Item(5).ColumnWidth = 41.6

# ---------------------------------------------------------------------
# 3. Checklist sheet: a handful of wrapped rows grew taller once Excel
#    recalculated the wrap height under the newer build.
# ---------------------------------------------------------------------
$checklist = $wb.Worksheets.Item("Checklist")

foreach ($r in @(3, 4, 5, 10, 11, 12, 17, 21, 26, 27)) {
    $checklist.Rows.Item($r).RowHeight = 30
}
foreach ($r in @(8, 20)) {
    $checklist.Rows.Item($r).RowHeight = 45
}

# ---------------------------------------------------------------------
# 4. Leave the selection where the author left it when saving
# ---------------------------------------------------------------------
$bvt.Activate()
$bvt.Range("E9").Select()
